# Generate Report for Handback
#
# 1) Update the Status text everywhere it is used (Overview!B2:C3 and the
#    "Status" column (C) on the zh-cn / de-de sheets all share the same
#    string).
# 2) Fill in the real handback timestamp that replaces the zero-date
#    placeholder in the "Latest Handback DateTime" column (H) on zh-cn,
#    and a newer handback timestamp on de-de.
# 3) Populate the previously-empty "Latest Target File" (F) and
#    "Latest Handback File" (G) columns on the zh-cn / de-de sheets with
#    hyperlinked file names, for both data rows.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$mdFile1 = "34dbf163-6549-49b7-ba19-b0521407b433.md"
$mdFile2 = "ffff0d5ceebf-7ee5-4e42-9be4-cb7241b79866.md"

$mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/2e4ee735c32e3df154f1abbc9c5700023bbe23ef/e2e/34dbf163-6549-49b7-ba19-b0521407b433.md"
$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/2e4ee735c32e3df154f1abbc9c5700023bbe23ef/e2e/ffff0d5ceebf-7ee5-4e42-9be4-cb7241b79866.md"

# ---- Overview sheet : Status text -------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

# ---- zh-cn sheet --------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zhXlf = "34dbf163-6549-49b7-ba19-b0521407b433.d20f0456541e02c2c6c4d06651f9e29ab6be105c.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/83def32953d9eb87ae05b7eea6e1ea0b96e52fab/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/34dbf163-6549-49b7-ba19-b0521407b433.d20f0456541e02c2c6c4d06651f9e29ab6be105c.zh-cn.xlf"

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus
$zh.Range("H2").Value = "2016-03-14 02:41:39"
$zh.Range("H3").Value = "2016-03-14 02:41:39"

$zh.Hyperlinks.Add($zh.Range("F2"), $mdUrl1, "", "", $mdFile1) | Out-Null
$zh.Hyperlinks.Add($zh.Range("G2"), $zhXlfUrl, "", "", $zhXlf) | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), $mdUrl2, "", "", $mdFile1) | Out-Null
$zh.Hyperlinks.Add($zh.Range("G3"), $zhXlfUrl, "", "", $zhXlf) | Out-Null

# Match the look of the other hyperlinked cells (A/B/D columns use the
# workbook's "HyperLink" style: underlined, accent-blue font).
foreach ($addr in @("F2", "G2", "F3", "G3")) {
    $zh.Range($addr).Font.Underline = 2
    $zh.Range($addr).Font.Color = 15570276
}

# ---- de-de sheet ----------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$deXlf = "34dbf163-6549-49b7-ba19-b0521407b433.d20f0456541e02c2c6c4d06651f9e29ab6be105c.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8c48e4b6299fca1a7b438d4631369a89acfbf63e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/34dbf163-6549-49b7-ba19-b0521407b433.d20f0456541e02c2c6c4d06651f9e29ab6be105c.de-de.xlf"

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus
$de.Range("H2").Value = "2016-03-14 02:41:44"
$de.Range("H3").Value = "2016-03-14 02:41:44"

$de.Hyperlinks.Add($de.Range("F2"), $mdUrl1, "", "", $mdFile1) | Out-Null
$de.Hyperlinks.Add($de.Range("G2"), $deXlfUrl, "", "", $deXlf) | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), $mdUrl2, "", "", $mdFile1) | Out-Null
$de.Hyperlinks.Add($de.Range("G3"), $deXlfUrl, "", "", $deXlf) | Out-Null

foreach ($addr in @("F2", "G2", "F3", "G3")) {
    $de.Range($addr).Font.Underline = 2
    $de.Range($addr).Font.Color = 15570276
}
